$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.17119999999999
$ws.Range("A4").Value = -21.57170000000001
$ws.Range("D4").Value = -7.931300000000004
$ws.Range("E4").Value = 12.54180000000001
$ws.Range("D5").Value = -8.546500000000002
$ws.Range("A6").Value = -21.64100000000001
$ws.Range("A7").Value = -21.455
$ws.Range("D8").Value = -8.343199999999998
$ws.Range("E9").Value = 13.44680000000001
$ws.Range("E11").Value = 13.2148
$ws.Range("E14").Value = 12.34350000000001
$ws.Range("A16").Value = -19.99939999999999
$ws.Range("D16").Value = -8.250599999999999
$ws.Range("E18").Value = 12.5566
$ws.Range("A20").Value = -22.08650000000003
$ws.Range("D22").Value = -7.819499999999994
$ws.Range("E25").Value = 13.1179
